$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the formatting of the existing
# header cells (e.g. G1 - bold/centered/bordered header style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for the "Save" column on row 2.
$ws.Range("H2").Value = 0
